# "edições para trabalhar no caos"
#
# The questionnaire's header row is re-labelled from plain sequence numbers
# (1..9) to "P1".."P9", the header row is centered (horizontally + vertically)
# with an explicit black font color, and the placeholder answer "_" in column
# E of the data row is replaced with "NA".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 1: turn the numeric column headers B1:J1 into P1..P9
$headers = @("P1", "P2", "P3", "P4", "P5", "P6", "P7", "P8", "P9")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 2).Value = $headers[$i]
}

# Row 2: the placeholder "_" answer becomes "NA"
$ws.Range("E2").Value = "NA"

# Header row formatting: centered, explicit black font color
$headerRange = $ws.Range("A1:J1")
$headerRange.Font.Color = 0
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4108

# Make the header row the active selection (matches the saved view state)
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollColumn = 1
$win.ScrollRow = 1
$ws.Range("A1:J1").Select()
